# Insert a new price record at row 354 (Ciboulette / Vega Modelo de Temuco),
# pushing the existing rows 354-401 down to 355-402.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 354; this shifts rows 354:401 -> 355:402
# and extends the sheet dimension to A1:R402.
$ws.Rows.Item(354).Insert()

# Populate the newly inserted row 354 with the new weekly record.
$ws.Cells.Item(354, 1).Value  = 10
$ws.Cells.Item(354, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(354, 3).Value  = "La Araucanía"
$ws.Cells.Item(354, 4).Value  = 45124
$ws.Cells.Item(354, 5).Value  = 9
$ws.Cells.Item(354, 6).Value  = 100112039
$ws.Cells.Item(354, 7).Value  = "Ciboulette"
$ws.Cells.Item(354, 8).Value  = "Sin especificar"
$ws.Cells.Item(354, 9).Value  = "Primera"
$ws.Cells.Item(354, 10).Value = 30
$ws.Cells.Item(354, 11).Value = 7000
$ws.Cells.Item(354, 12).Value = 7000
$ws.Cells.Item(354, 13).Value = 7000
$ws.Cells.Item(354, 14).Value = "$/docena de atados"
$ws.Cells.Item(354, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(354, 16).Value = 2333
$ws.Cells.Item(354, 17).Value = 3
$ws.Cells.Item(354, 18).Value = "Hortaliza"
